$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.598221302032471
$ws.Range("B1").Value = 2.835617542266846
$ws.Range("C1").Value = 2.245687484741211
$ws.Range("D1").Value = 2.09670090675354
$ws.Range("E1").Value = 1.794364213943481
